$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1953   # was 1952
$ws.Range("F12").Value = 1767   # was 1756
$ws.Range("F27").Value = 1033   # was 1032
$ws.Range("F28").Value = 4587   # was 4584
$ws.Range("F29").Value = 108   # was 107
$ws.Range("F34").Value = 1247   # was 1246

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 31   # was 30
$ws.Range("F19").Value = 188   # was 186
$ws.Range("F23").Value = 3   # was 2

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 9598   # was 9597
$ws.Range("F9").Value = 3106   # was 3105
$ws.Range("F10").Value = 620   # was 619
$ws.Range("F11").Value = 896   # was 895
$ws.Range("F12").Value = 321   # was 320
$ws.Range("F13").Value = 41   # was 39
$ws.Range("F14").Value = 60   # was 58
$ws.Range("F15").Value = 12   # was 10

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 3106   # was 3105
$ws.Range("F7").Value = 620   # was 619
$ws.Range("F8").Value = 896   # was 895
$ws.Range("F9").Value = 1953   # was 1952
$ws.Range("F10").Value = 41   # was 39
$ws.Range("F11").Value = 41   # was 39
$ws.Range("F12").Value = 60   # was 58
$ws.Range("F13").Value = 60   # was 58
$ws.Range("F16").Value = 12   # was 10
$ws.Range("F17").Value = 1767   # was 1756
$ws.Range("F26").Value = 188   # was 186
$ws.Range("F35").Value = 1033   # was 1032
$ws.Range("F40").Value = 108   # was 107
$ws.Range("F50").Value = 1247   # was 1246
